$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AON")

# Row 7 - Non-operating Income/Expense
$ws.Range("B7").Value = -235000000.0
$ws.Range("D7").Value = 47000000.0
$ws.Range("E7").Value = 50000000.0
$ws.Range("F7").Value = -27000000.0

# Row 11 - Income after Tax
$ws.Range("B11").Value = 2159000000.0
$ws.Range("D11").Value = 1864000000.0
$ws.Range("E11").Value = 1812000000.0
$ws.Range("F11").Value = 1689000000.0

# Row 14 - Gross Margin
$ws.Range("B14").Value = 0.4633

# Row 15 - EBIT Margin
$ws.Range("B15").Value = 0.2631

# Row 16 - EBT margin
$ws.Range("B16").Value = 0.233

# Row 17 - Net Profit Margin
$ws.Range("B17").Value = 0.1855

# Row 18 - Free Cash Flow Margin
$ws.Range("B18").Value = 0.2823

# Row 19 - EBITDA
$ws.Range("B19").Value = 3349000000.0

# Row 21 - EPS (Diluted, from Cont. Ops)
$ws.Range("B21").Value = 9.1672
$ws.Range("D21").Value = 7.8531
$ws.Range("E21").Value = 7.6031
$ws.Range("F21").Value = 7.0431

# Row 22 - EPS (Diluted, from Disc. Ops)
$ws.Range("C22").Value = -0.0043

# Row 24 - EPS (Basic, from Discontinued Ops)
$ws.Range("C24").Value = -0.0043

# Row 25 - EPS (Basic, from Continuous Ops)
$ws.Range("B25").Value = 9.2064
$ws.Range("D25").Value = 7.8919
$ws.Range("E25").Value = 7.6519
$ws.Range("F25").Value = 7.0919

# Row 26 - Income from Continuous Operations
$ws.Range("B26").Value = 2159000000.0
$ws.Range("D26").Value = 1864000000.0
$ws.Range("E26").Value = 1812000000.0
$ws.Range("F26").Value = 1689000000.0

# Row 27 - Income from Discontinued Operations
$ws.Range("C27").Value = 0.0

# Row 29 - EBITDA Margin
$ws.Range("B29").Value = 0.0

# Row 30 - Operating Cash Flow Margin
$ws.Range("B30").Value = 0.292
